# Update scripts with new TPM-derived expression/specificity/edge-weight values
# for the Col18a1-Ptprs ligand-receptor pair table (rows 2-10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.5789666666666667
$ws.Cells.Item(2, 8).Value = 1.7369
$ws.Cells.Item(2, 9).Value = 0.01523705650035473
$ws.Cells.Item(2, 10).Value = 0.01523705650035472
$ws.Cells.Item(2, 13).Value = 2.425633666666667
$ws.Cells.Item(2, 14).Value = 7.276901000000001
$ws.Cells.Item(2, 15).Value = 0.0662600404061536
$ws.Cells.Item(2, 16).Value = 0.06626004040615362
$ws.Cells.Item(2, 17).Value = 1.404361038544445
$ws.Cells.Item(2, 18).Value = 12.6392493469
$ws.Cells.Item(2, 19).Value = 0.00100960797938435
$ws.Cells.Item(2, 20).Value = 0.00100960797938435

$ws.Cells.Item(3, 7).Value = 0.5789666666666667
$ws.Cells.Item(3, 8).Value = 1.7369
$ws.Cells.Item(3, 9).Value = 0.01523705650035473
$ws.Cells.Item(3, 10).Value = 0.01523705650035472
$ws.Cells.Item(3, 15).Value = 0.4234968256437875
$ws.Cells.Item(3, 16).Value = 0.4234968256437876
$ws.Cells.Item(3, 17).Value = 8.975884080900002
$ws.Cells.Item(3, 18).Value = 80.7829567281
$ws.Cells.Item(3, 19).Value = 0.006452845060055265
$ws.Cells.Item(3, 20).Value = 0.006452845060055265

$ws.Cells.Item(4, 7).Value = 0.5789666666666667
$ws.Cells.Item(4, 8).Value = 1.7369
$ws.Cells.Item(4, 9).Value = 0.01523705650035473
$ws.Cells.Item(4, 10).Value = 0.01523705650035472
$ws.Cells.Item(4, 13).Value = 18.67887366666667
$ws.Cells.Item(4, 14).Value = 56.03662100000001
$ws.Cells.Item(4, 15).Value = 0.5102431339500588
$ws.Cells.Item(4, 16).Value = 0.5102431339500588
$ws.Cells.Item(4, 17).Value = 10.81444522387778
$ws.Cells.Item(4, 18).Value = 97.33000701490002
$ws.Cells.Item(4, 19).Value = 0.007774603460915111
$ws.Cells.Item(4, 20).Value = 0.00777460346091511

$ws.Cells.Item(5, 9).Value = 0.6545086962501954
$ws.Cells.Item(5, 10).Value = 0.6545086962501954
$ws.Cells.Item(5, 13).Value = 2.425633666666667
$ws.Cells.Item(5, 14).Value = 7.276901000000001
$ws.Cells.Item(5, 15).Value = 0.0662600404061536
$ws.Cells.Item(5, 16).Value = 0.06626004040615362
$ws.Cells.Item(5, 17).Value = 60.32441452067177
$ws.Cells.Item(5, 18).Value = 542.919730686046
$ws.Cells.Item(5, 19).Value = 0.04336777265971686
$ws.Cells.Item(5, 20).Value = 0.04336777265971687

$ws.Cells.Item(6, 9).Value = 0.6545086962501954
$ws.Cells.Item(6, 10).Value = 0.6545086962501954
$ws.Cells.Item(6, 15).Value = 0.4234968256437875
$ws.Cells.Item(6, 16).Value = 0.4234968256437876
$ws.Cells.Item(6, 19).Value = 0.2771823552182117
$ws.Cells.Item(6, 20).Value = 0.2771823552182117

$ws.Cells.Item(7, 9).Value = 0.6545086962501954
$ws.Cells.Item(7, 10).Value = 0.6545086962501954
$ws.Cells.Item(7, 13).Value = 18.67887366666667
$ws.Cells.Item(7, 14).Value = 56.03662100000001
$ws.Cells.Item(7, 15).Value = 0.5102431339500588
$ws.Cells.Item(7, 16).Value = 0.5102431339500588
$ws.Cells.Item(7, 17).Value = 464.5351576916852
$ws.Cells.Item(7, 18).Value = 4180.816419225166
$ws.Cells.Item(7, 19).Value = 0.3339585683722668
$ws.Cells.Item(7, 20).Value = 0.3339585683722668

$ws.Cells.Item(8, 9).Value = 0.33025424724945
$ws.Cells.Item(8, 10).Value = 0.3302542472494499
$ws.Cells.Item(8, 13).Value = 2.425633666666667
$ws.Cells.Item(8, 14).Value = 7.276901000000001
$ws.Cells.Item(8, 15).Value = 0.0662600404061536
$ws.Cells.Item(8, 16).Value = 0.06626004040615362
$ws.Cells.Item(8, 17).Value = 30.43870039073189
$ws.Cells.Item(8, 18).Value = 273.948303516587
$ws.Cells.Item(8, 19).Value = 0.0218826597670524
$ws.Cells.Item(8, 20).Value = 0.0218826597670524

$ws.Cells.Item(9, 9).Value = 0.33025424724945
$ws.Cells.Item(9, 10).Value = 0.3302542472494499
$ws.Cells.Item(9, 15).Value = 0.4234968256437875
$ws.Cells.Item(9, 16).Value = 0.4234968256437876
$ws.Cells.Item(9, 19).Value = 0.1398616253655206
$ws.Cells.Item(9, 20).Value = 0.1398616253655206

$ws.Cells.Item(10, 9).Value = 0.33025424724945
$ws.Cells.Item(10, 10).Value = 0.3302542472494499
$ws.Cells.Item(10, 13).Value = 18.67887366666667
$ws.Cells.Item(10, 14).Value = 56.03662100000001
$ws.Cells.Item(10, 15).Value = 0.5102431339500588
$ws.Cells.Item(10, 16).Value = 0.5102431339500588
$ws.Cells.Item(10, 17).Value = 234.3967462973586
$ws.Cells.Item(10, 18).Value = 2109.570716676228
$ws.Cells.Item(10, 19).Value = 0.1685099621168769
$ws.Cells.Item(10, 20).Value = 0.1685099621168769
